$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34/35: Monero and NEARProtocol swapped positions
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"

# Price (column D) updates
$ws.Range("D2").Value = "57.160.28"
$ws.Range("D3").Value = "3.068.95"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "3.069.62"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Value = "3.592.70"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Value = "57.169.67"
$ws.Range("D18").Value = "3.065.43"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "348.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.498"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Value = "0.0₃0843"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0653"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.688"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "2.408.43"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "3.106.53"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0260"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.927"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.30"
$ws.Range("D51").Style = "Normal"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("E6").Value = "  -5.31%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("E9").Value = "  +5.37%  "
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("E19").Value = "  -4.46%  "
$ws.Range("E20").Value = "  -3.57%  "
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("E28").Value = "  -10.01%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -4.82%  "
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  -10.36%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -6.27%  "
$ws.Range("E37").Value = "  -4.27%  "
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("E41").Value = "  -6.39%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("E44").Value = "  +5.56%  "
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("E50").Value = "  -8.40%  "
$ws.Range("E51").Value = "  -6.83%  "
